$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blind 75")

# Row 16 - Sum of Two Integers: Topic DP -> Binary
$ws.Range("C16").Value = "Binary"

# Row 28 - Longest Common Subsequence: Topic Graph -> DP
$ws.Range("C28").Value = "DP"

# Row 29 - Clone Graph: Topic Graph -> DP
$ws.Range("C29").Value = "DP"

# Row 53 - Longest Repeating Character Replacement: mark completed, add notes + runtime
$ws.Range("E53").Value = "X"
$ws.Range("F53").Value = "Dict of freq of chars, keep l, r pointer. Move r constantly, add r char to freq, if freq-l+r+1>k move l otherwise don" + [char]0x2019 + "t, take max of ans and r-l+1."
$ws.Range("G53").Value = "O(N)"

# Update selection to F53
$ws.Range("F53").Select()
